# Update cached market/profit values across multiple worksheets
# (scheduled runner refresh of Universalis price data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 62512.07
$ws.Range("I62").Value = 101196.336
$ws.Range("J62").Value = 51961.816
$ws.Range("K62").Value = 101196.336
$ws.Range("L62").Value = 51961.816
$ws.Range("M62").Value = -100572.336
$ws.Range("N62").Value = -53209.816

# Row 65
$ws.Range("H65").Value = 62512.07
$ws.Range("I65").Value = 101196.336
$ws.Range("J65").Value = 51961.816
$ws.Range("K65").Value = 505981.68
$ws.Range("L65").Value = 259809.08
$ws.Range("M65").Value = -502861.68
$ws.Range("N65").Value = -266049.08

# Row 69
$ws.Range("H69").Value = 8309.25
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 8309.25
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 24927.75
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -26675.75

# Row 72
$ws.Range("H72").Value = 8309.25
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 8309.25
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 74783.25
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -83519.25

# Row 137
$ws.Range("H137").Value = 3342.2727
$ws.Range("I137").Value = 2798.5
$ws.Range("K137").Value = 8395.5
$ws.Range("M137").Value = -5845.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2507.7678
$ws.Range("I32").Value = 2470
$ws.Range("K32").Value = 2470
$ws.Range("M32").Value = -2183

# Row 45
$ws.Range("H45").Value = 1772
$ws.Range("J45").Value = 2099.6
$ws.Range("L45").Value = 2099.6
$ws.Range("N45").Value = -2853.6

# Row 61
$ws.Range("H61").Value = 66669612
$ws.Range("I61").Value = 83336210
$ws.Range("J61").Value = 3214
$ws.Range("K61").Value = 83336210
$ws.Range("L61").Value = 3214
$ws.Range("M61").Value = -83335998
$ws.Range("N61").Value = -3638

# Row 102
$ws.Range("H102").Value = 8816.52
$ws.Range("I102").Value = 2726.6843
$ws.Range("J102").Value = 28101
$ws.Range("K102").Value = 2726.6843
$ws.Range("L102").Value = 28101
$ws.Range("M102").Value = -1104.6843
$ws.Range("N102").Value = -31345

# Row 122
$ws.Range("H122").Value = 19611470
$ws.Range("I122").Value = 30305908
$ws.Range("K122").Value = 90917724
$ws.Range("M122").Value = -90915274

# Row 136
$ws.Range("H136").Value = 66669612
$ws.Range("I136").Value = 83336210
$ws.Range("J136").Value = 3214
$ws.Range("K136").Value = 250008630
$ws.Range("L136").Value = 9642
$ws.Range("M136").Value = -250006080
$ws.Range("N136").Value = -14742

$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 69000
$ws.Range("J100").Value = 69000
$ws.Range("L100").Value = 69000
$ws.Range("N100").Value = -71164

# Row 109
$ws.Range("H109").Value = 130000.5
$ws.Range("J109").Value = 130000.5
$ws.Range("L109").Value = 130000.5
$ws.Range("N109").Value = -132774.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2601.3809
$ws.Range("I31").Value = 1833.0714
$ws.Range("J31").Value = 4138
$ws.Range("K31").Value = 1833.0714
$ws.Range("L31").Value = 4138
$ws.Range("M31").Value = -1538.0714
$ws.Range("N31").Value = -4728

# Row 34
$ws.Range("H34").Value = 2601.3809
$ws.Range("I34").Value = 1833.0714
$ws.Range("J34").Value = 4138
$ws.Range("K34").Value = 1833.0714
$ws.Range("L34").Value = 4138
$ws.Range("M34").Value = -1631.0714
$ws.Range("N34").Value = -4542

# Row 74
$ws.Range("H74").Value = 33676.75
$ws.Range("J74").Value = 33676.75
$ws.Range("L74").Value = 33676.75
$ws.Range("N74").Value = -35424.75

# Row 77
$ws.Range("H77").Value = 33676.75
$ws.Range("J77").Value = 33676.75
$ws.Range("L77").Value = 101030.25
$ws.Range("N77").Value = -109766.25

# Row 122
$ws.Range("H122").Value = 1766.1052
$ws.Range("I122").Value = 1766.1052
$ws.Range("K122").Value = 5298.3156
$ws.Range("M122").Value = -2848.3156

# Row 134
$ws.Range("H134").Value = 2708.889
$ws.Range("I134").Value = 2297.625
$ws.Range("K134").Value = 6892.875
$ws.Range("M134").Value = -4357.875

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 302
$ws.Range("I15").Value = 107.85714
$ws.Range("J15").Value = 755
$ws.Range("K15").Value = 323.57142
$ws.Range("L15").Value = 2265
$ws.Range("M15").Value = -183.57142
$ws.Range("N15").Value = -2545

# Row 102
$ws.Range("H102").Value = 26250
$ws.Range("I102").Value = 2500
$ws.Range("J102").Value = 50000
$ws.Range("K102").Value = 7500
$ws.Range("L102").Value = 150000
$ws.Range("M102").Value = -5066
$ws.Range("N102").Value = -154868

# Row 132
$ws.Range("H132").Value = 2086
$ws.Range("I132").Value = 1425.25
$ws.Range("J132").Value = 2614.6
$ws.Range("K132").Value = 12827.25
$ws.Range("L132").Value = 23531.4
$ws.Range("M132").Value = -10297.25
$ws.Range("N132").Value = -28591.4

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1549.7354
$ws.Range("I102").Value = 761.75
$ws.Range("K102").Value = 761.75
$ws.Range("M102").Value = 860.25

# Row 122
$ws.Range("H122").Value = 15154048
$ws.Range("I122").Value = 1734.1818
$ws.Range("J122").Value = 45458676
$ws.Range("K122").Value = 5202.5454
$ws.Range("L122").Value = 136376028
$ws.Range("M122").Value = -2752.5454
$ws.Range("N122").Value = -136380928

# Row 132
$ws.Range("H132").Value = 3413.5293
$ws.Range("I132").Value = 3008.423
$ws.Range("K132").Value = 9025.269
$ws.Range("M132").Value = -6495.269

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1900.2
$ws.Range("I40").Value = 1808.5834
$ws.Range("J40").Value = 2266.6667
$ws.Range("K40").Value = 1808.5834
$ws.Range("L40").Value = 2266.6667
$ws.Range("M40").Value = -1672.5834
$ws.Range("N40").Value = -2538.6667

# Row 93
$ws.Range("H93").Value = 6626.25
$ws.Range("I93").Value = 6833.6665
$ws.Range("J93").Value = 6004
$ws.Range("K93").Value = 6833.6665
$ws.Range("L93").Value = 6004
$ws.Range("M93").Value = -5585.6665
$ws.Range("N93").Value = -8500

# Row 100
$ws.Range("H100").Value = 2427.3333
$ws.Range("I100").Value = 2149.375
$ws.Range("K100").Value = 2149.375
$ws.Range("M100").Value = -1608.375

# Row 122
$ws.Range("H122").Value = 3448.4583
$ws.Range("I122").Value = 2709.0557
$ws.Range("K122").Value = 8127.1671
$ws.Range("M122").Value = -5677.1671

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 60070
$ws.Range("I51").Value = 60070
$ws.Range("K51").Value = 60070
$ws.Range("M51").Value = -59560

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()

# Row 81
$ws.Range("H81").Value = 9532041
$ws.Range("I81").Value = 6560.2856
$ws.Range("J81").Value = 14294781
$ws.Range("K81").Value = 13120.5712
$ws.Range("L81").Value = 28589562
$ws.Range("M81").Value = -12059.5712
$ws.Range("N81").Value = -28591684

# Row 84
$ws.Range("H84").Value = 9532041
$ws.Range("I84").Value = 6560.2856
$ws.Range("J84").Value = 14294781
$ws.Range("K84").Value = 65602.856
$ws.Range("L84").Value = 142947810
$ws.Range("M84").Value = -60298.856
$ws.Range("N84").Value = -142958418

# Row 95
$ws.Range("H95").Value = 33483.832
$ws.Range("J95").Value = 33483.832
$ws.Range("L95").Value = 33483.832
$ws.Range("N95").Value = -38975.832

# Row 96
$ws.Range("H96").Value = 5198.6665
$ws.Range("I96").Value = 3557
$ws.Range("J96").Value = 7250.75
$ws.Range("K96").Value = 3557
$ws.Range("L96").Value = 7250.75
$ws.Range("M96").Value = -2184
$ws.Range("N96").Value = -9996.75

# Row 126
$ws.Range("H126").Value = 1583.6316
$ws.Range("I126").Value = 1440.5294
$ws.Range("K126").Value = 4321.5882
$ws.Range("M126").Value = -1851.5882

# Row 132
$ws.Range("H132").Value = 4769.96
$ws.Range("I132").Value = 4929.4116
$ws.Range("K132").Value = 14788.2348
$ws.Range("M132").Value = -12258.2348

Write-Host "Updated profit sheets with latest Universalis data"